$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet: Date value ---
$wsMeta.Range("B8").Value = "2025-05-05T14:17:01+00:00"

# --- Elements sheet ---

# Binding Strength: required -> preferred (shared string used by every row below,
# so every "required" cell in column X becomes "preferred")
$wsElem.Range("X6").Value = "preferred"
$wsElem.Range("X22").Value = "preferred"
$wsElem.Range("X23").Value = "preferred"
$wsElem.Range("X25").Value = "preferred"
$wsElem.Range("X26").Value = "preferred"
$wsElem.Range("X27").Value = "preferred"
$wsElem.Range("X28").Value = "preferred"
$wsElem.Range("X35").Value = "preferred"

# legalAuthenticator row (row 11): path renamed to legalAuthenticator[x]
$wsElem.Range("A11").Value = "DocumentEntry.legalAuthenticator[x]"
$wsElem.Range("B11").Value = "DocumentEntry.legalAuthenticator[x]"
$wsElem.Range("AF11").Value = "DocumentEntry.legalAuthenticator[x]"

# Type(s) column for that row: drop the Reference(...) wrapper syntax
$wsElem.Range("K11").Value = "https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorPS" + [char]10 + "https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorPatienthttps://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/ActorSystem"

# Binding Value Set URLs: DMP nomenclature -> CISIS nomenclature
$wsElem.Range("Z22").Value = "https://mos.esante.gouv.fr/NOS/JDV_J06-XdsClassCode-CISIS/FHIR/JDV-J06-XdsClassCode-CISIS"
$wsElem.Range("Z25").Value = "https://mos.esante.gouv.fr/NOS/JDV_J10-XdsFormatCode-CISIS/FHIR/JDV-J10-XdsFormatCode-CISIS"
$wsElem.Range("Z26").Value = "https://mos.esante.gouv.fr/NOS/JDV_J02-XdsHealthcareFacilityTypeCode-CISIS/FHIR/JDV-J02-XdsHealthcareFacilityTypeCode-CISIS"
$wsElem.Range("Z27").Value = "https://mos.esante.gouv.fr/NOS/JDV_J04-XdsPracticeSettingCode-CISIS/FHIR/JDV-J04-XdsPracticeSettingCode-CISIS"
$wsElem.Range("Z28").Value = "https://mos.esante.gouv.fr/NOS/JDV_J07-XdsTypeCode-CISIS/FHIR/JDV-J07-XdsTypeCode-CISIS"
$wsElem.Range("Z35").Value = "https://mos.esante.gouv.fr/NOS/JDV_J197-XdsTypesIdentifiantsReferenceId-CISIS/FHIR/JDV-J197-XdsTypesIdentifiantsReferenceId-CISIS"

# Column widths follow content auto-fit (K shrank, Z grew slightly) on the Elements sheet
$wsElem.Columns.Item(11).ColumnWidth = 128.05078125
$wsElem.Columns.Item(26).ColumnWidth = 106.83203125
